$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Goldfish" entry (row 41) was removed; every row below it shifts up
# by one.
$ws.Rows.Item(41).Delete()

# A new "South African abalone" entry was added right before "Southern
# rock lobster" (which, after the deletion above, now sits at row 80).
$ws.Rows.Item(80).Insert()
$ws.Cells.Item(80, 1).Value = "South African abalone"
$ws.Cells.Item(80, 2).Value = "Haliotis midae"
$ws.Cells.Item(80, 3).Value = "invert"
$ws.Cells.Item(80, 4).Value = "Gastropoda"
$ws.Cells.Item(80, 5).Value = "Patellogastropoda"
$ws.Cells.Item(80, 6).Value = "Haliotidae"
$ws.Cells.Item(80, 7).Value = "Haliotis"

# Fix a typo: "Tebula snail" should read "Tegula snail" (matches the
# existing sci_name "Tegula spp." already in that row).
$ws.Cells.Item(86, 1).Value = "Tegula snail"
